# Update "Execution Flag" column (E) for the UITestCases sheet.
# Rows 2-11 and 13-28 switch from "Yes" to "No"; row 12 is left as "Yes".
# Rows 29-34 were already "No" and stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UITestCases")

$rowsToMarkNo = @(2,3,4,5,6,7,8,9,10,11,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28)

foreach ($r in $rowsToMarkNo) {
    $ws.Range("E$r").Value = "No"
}

# Reflect the user's final selection/cursor position on the sheet.
$null = $ws.Range("E12").Select()
